# Edit: add "PO Forecast" sheet, rename headers on existing sheets.

$wb = $excel.ActiveWorkbook

# 1) Rename header cells on the two existing sheets.
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 2) Add the new "PO Forecast" sheet after "Monthly Trend".
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Reuse the existing header formatting (bold, centered, bordered) from the
# "Weekly Quantity" sheet's header row, and the date-like number formatting
# from its "A" data column, so the new sheet matches existing look & feel.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A30").PasteSpecial(-4122)

# 3) Header row values.
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# 4) Data rows: ds, PO_Forecast, yhat_lower, yhat_upper
$rows = @(
    @(2, 45312.99999999999, 242, -94.13812896678213, 590.3573777675928),
    @(3, 45319.99999999999, 243, -107.6060260361539, 563.8462540693727),
    @(4, 45354.99999999999, 251, -93.25275333630549, 581.7417407432782),
    @(5, 45375.99999999999, 256, -71.57863208608477, 578.129330786212),
    @(6, 45389.99999999999, 259, -63.09985243869867, 601.5031032211979),
    @(7, 45410.99999999999, 264, -70.44423504369551, 612.8575336947431),
    @(8, 45417.99999999999, 265, -59.01297789710681, 585.7434500551666),
    @(9, 45431.99999999999, 268, -55.13162285470388, 624.7157979407846),
    @(10, 45438.99999999999, 270, -29.15765355612379, 631.7208723334394),
    @(11, 45445.99999999999, 271, -52.04740321654521, 586.9524522814169),
    @(12, 45459.99999999999, 275, -46.43187421019154, 591.638197111783),
    @(13, 45466.99999999999, 276, -41.77036743650535, 622.9811161853366),
    @(14, 45473.99999999999, 278, -57.77381991664202, 582.9950899741617),
    @(15, 45529.99999999999, 290, -44.42366055798218, 603.8525456212676),
    @(16, 45536.99999999999, 292, -34.46718662383869, 613.4969506922179),
    @(17, 45543.99999999999, 294, -17.65331028280268, 628.7010365471905),
    @(18, 45550.99999999999, 295, -26.05378888747623, 619.1989175442332),
    @(19, 45571.99999999999, 300, 1.982110082768023, 649.0012117036763),
    @(20, 45585.99999999999, 303, -44.11115102949659, 637.3289152226879),
    @(21, 45592.99999999999, 305, -45.10894708639017, 637.1132637631767),
    @(22, 45599.99999999999, 306, -6.959574240167655, 669.3855109319617),
    @(23, 45606.99999999999, 308, -35.25915263275606, 627.8987699662498),
    @(24, 45613.99999999999, 309, -6.807488529380513, 640.3316011723504),
    @(25, 45620.99999999999, 311, -21.38233839027835, 675.1195579822972),
    @(26, 45627.99999999999, 312, -24.13833038917155, 651.8755150613601),
    @(27, 45634.99999999999, 314, -28.57927132787847, 651.1979535647362),
    @(28, 45641.99999999999, 316, -47.07527908462793, 627.8941506169899),
    @(29, 45648.99999999999, 317, 5.910626759465179, 644.4110066706933),
    @(30, 45655.99999999999, 319, -13.45032502940971, 647.9615538566065)
)

foreach ($row in $rows) {
    $r = $row[0]
    $wsForecast.Cells.Item($r, 1).Value = $row[1]
    $wsForecast.Cells.Item($r, 2).Value = $row[2]
    $wsForecast.Cells.Item($r, 3).Value = $row[3]
    $wsForecast.Cells.Item($r, 4).Value = $row[4]
}
